$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/10_tokiti1.wav"
$ws.Range("B2").Value = "pngimages/10_backpack.png"

$ws.Range("A3").Value = "trainingaudio/19_papipi1.wav"
$ws.Range("B3").Value = "pngimages/19_burger.png"

$ws.Range("A4").Value = "trainingaudio/21_papika1.wav"
$ws.Range("B4").Value = "pngimages/21_cheese.png"

$ws.Range("A5").Value = "trainingaudio/14_pokoto1.wav"
$ws.Range("B5").Value = "pngimages/14_coffee.png"

$ws.Range("A6").Value = "trainingaudio/09_tipata2.wav"
$ws.Range("B6").Value = "pngimages/09_plane.png"

$ws.Range("A7").Value = "trainingaudio/22_kakoki1.wav"
$ws.Range("B7").Value = "pngimages/22_egg.png"
